# Apply corrections to Summary, Repayment schedule and Transactions sheets.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("F2").Value = 0
# Touch G2 so a (blank, default-styled) cell gets materialised in the sheet,
# matching the new <c r="G2"/> entry in the target workbook.
$wsSummary.Range("G2").Style = "Normal"

$wsSummary.Range("A3").Value = 720.4
$wsSummary.Range("E3").Value = 520.4

$wsSummary.Range("D5").Select()

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Remove the now unused P2 and O3:O13 cells entirely (clear format first,
# then contents, so the cell node itself disappears from the XML).
$wsSchedule.Range("P2").ClearFormats()
$wsSchedule.Range("P2").ClearContents()

$wsSchedule.Range("O3:O13").ClearFormats()
$wsSchedule.Range("O3:O13").ClearContents()

# Row 4
$wsSchedule.Range("B4").Value = 31
$wsSchedule.Range("C4").Value = 42095
$wsSchedule.Range("F4").Value = 872.19
$wsSchedule.Range("G4").Value = 8363.27
$wsSchedule.Range("H4").Value = 92.35

# Row 5
$wsSchedule.Range("B5").Value = 30
$wsSchedule.Range("C5").Value = 42125
$wsSchedule.Range("F5").Value = 880.91
$wsSchedule.Range("G5").Value = 7482.36
$wsSchedule.Range("H5").Value = 83.63

# Row 6
$wsSchedule.Range("B6").Value = 31
$wsSchedule.Range("C6").Value = 42156
$wsSchedule.Range("F6").Value = 889.72
$wsSchedule.Range("G6").Value = 6592.64
$wsSchedule.Range("H6").Value = 74.819999999999993

# Row 7
$wsSchedule.Range("B7").Value = 30
$wsSchedule.Range("C7").Value = 42186
$wsSchedule.Range("F7").Value = 898.61
$wsSchedule.Range("G7").Value = 5694.03
$wsSchedule.Range("H7").Value = 65.930000000000007

# Row 8
$wsSchedule.Range("B8").Value = 31
$wsSchedule.Range("C8").Value = 42217
$wsSchedule.Range("F8").Value = 907.6
$wsSchedule.Range("G8").Value = 4786.43
$wsSchedule.Range("H8").Value = 56.94

# Row 9 (Days stays 31)
$wsSchedule.Range("C9").Value = 42248
$wsSchedule.Range("F9").Value = 916.68
$wsSchedule.Range("G9").Value = 3869.75
$wsSchedule.Range("H9").Value = 47.86

# Row 10
$wsSchedule.Range("B10").Value = 30
$wsSchedule.Range("C10").Value = 42278
$wsSchedule.Range("F10").Value = 925.84
$wsSchedule.Range("G10").Value = 2943.91
$wsSchedule.Range("H10").Value = 38.700000000000003

# Row 11
$wsSchedule.Range("B11").Value = 31
$wsSchedule.Range("C11").Value = 42309
$wsSchedule.Range("F11").Value = 935.1
$wsSchedule.Range("G11").Value = 2008.81
$wsSchedule.Range("H11").Value = 29.44

# Row 12 - G12 switches from the plain style to the currency style (like G4:G11)
$wsSchedule.Range("B12").Value = 30
$wsSchedule.Range("C12").Value = 42339
$wsSchedule.Range("F12").Value = 944.45
$wsSchedule.Range("G12").NumberFormat = $wsSchedule.Range("G11").NumberFormat
$wsSchedule.Range("G12").Value = 1064.3599999999999
$wsSchedule.Range("H12").Value = 20.09

# Row 13 - F13 switches to the currency style (like G2:G12); K13 and P13
# switch to the integer/thousands style used elsewhere (e.g. G2).
$wsSchedule.Range("B13").Value = 31
$wsSchedule.Range("C13").Value = 42370
$wsSchedule.Range("F13").NumberFormat = $wsSchedule.Range("G11").NumberFormat
$wsSchedule.Range("F13").Value = 1064.3599999999999
$wsSchedule.Range("H13").Value = 10.64
$wsSchedule.Range("K13").NumberFormat = $wsSchedule.Range("G2").NumberFormat
$wsSchedule.Range("K13").Value = 1075
$wsSchedule.Range("P13").NumberFormat = $wsSchedule.Range("G2").NumberFormat
$wsSchedule.Range("P13").Value = 1075

$wsSchedule.Range("H4:H13").Select()

# ---------------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")

$wsTransactions.Range("A2").Value = 6356
$wsTransactions.Range("A3").Value = 6354

# Select last so this sheet remains the active tab, as in the original file.
$wsTransactions.Range("D3").Select()
